# Split the "Lineage / ProductStrain" paragraph into two paragraphs and
# drop the stray <w:proofErr> markers, merging the "{" + "{Label1.ProductStrain}}"
# runs into a single "{{Label1.ProductStrain}}" run that lives on the new paragraph.

$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("{{Label1.Lineage}}{{Label1.ProductStrain}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the Lineage/ProductStrain paragraph"
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="22C2A3DD" w14:textId="2EDC11CD" w:rsidR="00E80EEC" w:rsidRPr="005A7DEA" w:rsidRDefault="00200931" w:rsidP="009D2E74"><w:pPr><w:spacing w:line="216" w:lineRule="auto"/><w:ind w:right="126"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>{{</w:t></w:r><w:r w:rsidR="00053E0C" w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Label1</w:t></w:r><w:r w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00053E0C" w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Lineage</w:t></w:r><w:r w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="216" w:lineRule="auto"/><w:ind w:right="126"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidR="002F4FA8" w:rsidRPr="008F15C5"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="8"/><w:szCs w:val="8"/></w:rPr><w:t>{{Label1.ProductStrain}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

Write-Output "Split complete. Paragraph count is now $($d.Paragraphs.Count)."
